$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set cell values to add "N" (north door) suffix to the existing room-letter codes
# (order matches the shared-strings table insertion order in the target file)
$ws.Range("B2").Value = "AN"
$ws.Range("H4").Value = "SN"
$ws.Range("N2").Value = "CN"
$ws.Range("S9").Value = "BN"
$ws.Range("Q13").Value = "MN"
$ws.Range("B12").Value = "EN"
$ws.Range("B20").Value = "TN"
$ws.Range("J20").Value = "RN"
$ws.Range("R22").Value = "PN"

# Update the selection on the sheet view
$ws.Range("R19").Select()
